{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// Applies the textual edits described by the commit\n// \"f\u00e9lre\u00e9rthet\u0151s\u00e9g \u00e9s helyes\u00edr\u00e1s jav\u00edt\u00e1sa\" (fixing ambiguity / spelling):\n//\n//  1. Rewrite the \"Ha nincs felhaszn\u00e1l\u00f3 fi\u00f3kja...\" sentence about the\n//     registration button.\n//  2. Drop a stray double space before \"A felhaszn\u00e1l\u00f3i fi\u00f3k t\u00f6rl\u00e9se...\"\n//     sentence (remove the trailing extra space run at paragraph end).\n//  3. Clarify \"Nagyobb k\u00e9peket helyez\u00fcnk el a slider-ben.\" ->\n//     \"Az \u00e9tteremr\u0151l k\u00e9peket helyez\u00fcnk el a slider-ben.\"\n//  4. \"Itt lehet a hibajelent\u00e9sre gombra kapcsolni.\" ->\n//     \"Itt lehet a hibajelent\u00e9s gombra kapcsolni.\"\n//  5. Add a missing comma: \"...jelenik meg ha be vannak...\" ->\n//     \"...jelenik meg, ha be vannak...\"\n//  6. Insert a new paragraph clarifying the minimum reservation notice\n//     period, right after \"Lesz lehet\u0151s\u00e9g a helyfoglal\u00e1sra.\"\n\nasync function replaceOnce(context, searchText, replacement) {\n  const results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. \"Ha nincs felhaszn\u00e1l\u00f3 fi\u00f3kja...\" -> \"Van egy regisztr\u00e1ci\u00f3 gomb ami \u00e1tviszi a regisztr\u00e1ci\u00f3s oldalra.\"\nawait replaceOnce(\n  context,\n  \"Ha nincs felhaszn\u00e1l\u00f3 fi\u00f3kja, akkor alatta lesz egy regisztr\u00e1ci\u00f3 gomb, ami elviszi a felhaszn\u00e1l\u00f3t a regisztr\u00e1ci\u00f3 oldalra.\",\n  \"Van egy regisztr\u00e1ci\u00f3 gomb ami \u00e1tviszi a regisztr\u00e1ci\u00f3s oldalra.\"\n);\n\n// 2. Remove the stray trailing space after \"...kattint\u00e1sra t\u00f6rt\u00e9nik. \"\nawait replaceOnce(\n  context,\n  \"A felhaszn\u00e1l\u00f3i fi\u00f3k t\u00f6rl\u00e9se egy gomb kattint\u00e1sra t\u00f6rt\u00e9nik. \",\n  \"A felhaszn\u00e1l\u00f3i fi\u00f3k t\u00f6rl\u00e9se egy gomb kattint\u00e1sra t\u00f6rt\u00e9nik.\"\n);\n\n// 3. \"Nagyobb k\u00e9peket helyez\u00fcnk el a slider-ben. \" -> \"Az \u00e9tteremr\u0151l k\u00e9peket helyez\u00fcnk el a slider-ben.\"\nawait replaceOnce(\n  context,\n  \"Nagyobb k\u00e9peket helyez\u00fcnk el a slider-ben. \",\n  \"Az \u00e9tteremr\u0151l k\u00e9peket helyez\u00fcnk el a slider-ben.\"\n);\n\n// 4. \" Itt lehet a hibajelent\u00e9sre\" -> \" Itt lehet a hibajelent\u00e9s\"\nawait replaceOnce(\n  context,\n  \" Itt lehet a hibajelent\u00e9sre\",\n  \" Itt lehet a hibajelent\u00e9s\"\n);\n\n// 5. \"Ami csak akkor jelenik meg ha be vannak jelentkezve.\" -> add comma\nawait replaceOnce(\n  context,\n  \"Ami csak akkor jelenik meg ha be vannak jelentkezve.\",\n  \"Ami csak akkor jelenik meg, ha be vannak jelentkezve.\"\n);\n\n// 6. Insert a new paragraph after \"Lesz lehet\u0151s\u00e9g a helyfoglal\u00e1sra.\"\nconst anchor = context.document.body.search(\"Lesz lehet\u0151s\u00e9g a helyfoglal\u00e1sra.\", { matchCase: true });\nanchor.load(\"items\");\nawait context.sync();\nif (anchor.items.length === 0) {\n  throw new Error(\"Anchor paragraph not found for insertion.\");\n}\nconst anchorParagraph = anchor.items[0].paragraphs.getFirst();\nanchorParagraph.insertParagraph(\n  \"\\tA foglal\u00e1s minimum 12 \u00f3r\u00e1val el\u0151re kell jelezni. \",\n  Word.InsertLocation.after\n);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word / $d (ActiveDocument) are pre-seeded by the harness.\n#\n# Applies the textual edits described by the commit\n# \"f\u00e9lre\u00e9rthet\u0151s\u00e9g \u00e9s helyes\u00edr\u00e1s jav\u00edt\u00e1sa\" (fixing ambiguity / spelling):\n#\n#  1. Rewrite the \"Ha nincs felhaszn\u00e1l\u00f3 fi\u00f3kja...\" sentence about the\n#     registration button.\n#  2. Drop a stray double space before \"A felhaszn\u00e1l\u00f3i fi\u00f3k t\u00f6rl\u00e9se...\"\n#     sentence (remove the trailing extra space run at paragraph end).\n#  3. Clarify \"Nagyobb k\u00e9peket helyez\u00fcnk el a slider-ben.\" ->\n#     \"Az \u00e9tteremr\u0151l k\u00e9peket helyez\u00fcnk el a slider-ben.\"\n#  4. \"Itt lehet a hibajelent\u00e9sre gombra kapcsolni.\" ->\n#     \"Itt lehet a hibajelent\u00e9s gombra kapcsolni.\"\n#  5. Add a missing comma: \"...jelenik meg ha be vannak...\" ->\n#     \"...jelenik meg, ha be vannak...\"\n#  6. Insert a new paragraph clarifying the minimum reservation notice\n#     period, right after \"Lesz lehet\u0151s\u00e9g a helyfoglal\u00e1sra.\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once($findText, $replaceText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $result) {\n        Write-Output \"NOT FOUND: $findText\"\n    }\n}\n\n# 1. \"Ha nincs felhaszn\u00e1l\u00f3 fi\u00f3kja...\" -> \"Van egy regisztr\u00e1ci\u00f3 gomb ami \u00e1tviszi a regisztr\u00e1ci\u00f3s oldalra.\"\nReplace-Once \"Ha nincs felhaszn\u00e1l\u00f3 fi\u00f3kja, akkor alatta lesz egy regisztr\u00e1ci\u00f3 gomb, ami elviszi a felhaszn\u00e1l\u00f3t a regisztr\u00e1ci\u00f3 oldalra.\" \"Van egy regisztr\u00e1ci\u00f3 gomb ami \u00e1tviszi a regisztr\u00e1ci\u00f3s oldalra.\"\n\n# 2. Remove the stray trailing space after \"...kattint\u00e1sra t\u00f6rt\u00e9nik. \"\nReplace-Once \"A felhaszn\u00e1l\u00f3i fi\u00f3k t\u00f6rl\u00e9se egy gomb kattint\u00e1sra t\u00f6rt\u00e9nik. \" \"A felhaszn\u00e1l\u00f3i fi\u00f3k t\u00f6rl\u00e9se egy gomb kattint\u00e1sra t\u00f6rt\u00e9nik.\"\n\n# 3. \"Nagyobb k\u00e9peket helyez\u00fcnk el a slider-ben. \" -> \"Az \u00e9tteremr\u0151l k\u00e9peket helyez\u00fcnk el a slider-ben.\"\nReplace-Once \"Nagyobb k\u00e9peket helyez\u00fcnk el a slider-ben. \" \"Az \u00e9tteremr\u0151l k\u00e9peket helyez\u00fcnk el a slider-ben.\"\n\n# 4. \" Itt lehet a hibajelent\u00e9sre\" -> \" Itt lehet a hibajelent\u00e9s\"\nReplace-Once \" Itt lehet a hibajelent\u00e9sre\" \" Itt lehet a hibajelent\u00e9s\"\n\n# 5. \"Ami csak akkor jelenik meg ha be vannak jelentkezve.\" -> add comma\nReplace-Once \"Ami csak akkor jelenik meg ha be vannak jelentkezve.\" \"Ami csak akkor jelenik meg, ha be vannak jelentkezve.\"\n\n# 6. Insert a new paragraph after \"Lesz lehet\u0151s\u00e9g a helyfoglal\u00e1sra.\"\n$paras = $d.Paragraphs\n$count = $paras.Count\n$targetIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Range.Text -like \"Lesz lehet\u0151s\u00e9g a helyfoglal\u00e1sra.*\") {\n        $targetIndex = $i\n        break\n    }\n}\nif ($targetIndex -gt 0) {\n    $p = $paras.Item($targetIndex)\n    $p.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Item($targetIndex + 1)\n    $newPara.Range.Text = \"`tA foglal\u00e1s minimum 12 \u00f3r\u00e1val el\u0151re kell jelezni. \"\n} else {\n    Write-Output \"NOT FOUND: anchor paragraph for insertion\"\n}\n"}
